$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing data rows (old rows 2-21) down to rows 3-22
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (the header); clear it so the
# new data row matches the plain (unstyled) formatting of the other data rows.
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new data point
$newRow2 = @(-0.002831595284598272, -0.19119707601411, 0.1513774214046346)
$ws.Cells.Item(2, 1).Value = $newRow2[0]
$ws.Cells.Item(2, 2).Value = $newRow2[1]
$ws.Cells.Item(2, 3).Value = $newRow2[2]

# Append nine new data rows after the existing data (now ending at row 22), i.e. rows 23-31
$appendData = @(
    @(-2.750307172536849, -1.468972404088292, 0.1615269269261984),
    @(-3.695782780647288, -0.6018874943256158, 1.057165026664715),
    @(-5.04273155757359, 2.776979684829711, -3.505371774945944),
    @(-2.32667221341812, 0.9583009992326832, -5.241142443248169),
    @(2.433021928582876, -2.759944068534038, 3.60468020609447),
    @(-2.125196490968979, 0.5359780830996369, 2.942476987838745),
    @(-3.455752406801498, 0.2274821900895648, 2.951905420848299),
    @(-2.549259322030204, 0.480571014540536, 3.463431903294155),
    @(0.6998523473739624, -1.159572852775452, -0.4669593572616656)
)

$startRow = 23
for ($i = 0; $i -lt $appendData.Count; $i++) {
    $r = $startRow + $i
    $row = $appendData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
